$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the previous 6 data rows (2-7); new data will repopulate rows 2-10
$ws.Range("A2:T7").ClearContents()

$newRows = @(
    @{A="ECs"; B="Lama2"; C="Itga7"; D="ECs"; E=2; F=0.6666666666666666; G=2.452389333333333; H=7.357168; I=0.007993767302975028; J=0.007993767302975028; K=3; L=1; M=3.349417; N=10.048251; O=0.0602955800561437; P=0.0602955800561437; Q=8.214074523685333; R=73.926670713168; S=0.0004819888363667147; T=0.0004819888363667146},
    @{A="ECs"; B="Lama2"; C="Itga7"; D="FAPs"; E=2; F=0.6666666666666666; G=2.452389333333333; H=7.357168; I=0.007993767302975028; J=0.007993767302975028; K=3; L=1; M=1.282876; N=3.848628; O=0.023094094452887; P=0.02309409445288699; Q=3.146111418389333; R=28.315002765504; S=0.000184608817129305; T=0.000184608817129305},
    @{A="ECs"; B="Lama2"; C="Itga7"; D="sCs"; E=2; F=0.6666666666666666; G=2.452389333333333; H=7.357168; I=0.007993767302975028; J=0.007993767302975028; K=3; L=1; M=50.917666; N=152.752998; O=0.9166103254909692; P=0.9166103254909692; Q=124.8699409766293; R=1123.829468789664; S=0.007327169649479007; T=0.007327169649479007},
    @{A="FAPs"; B="Lama2"; C="Itga7"; D="ECs"; E=3; F=1; G=243.5672963333334; H=730.701889; I=0.7939278902575405; J=0.7939278902575405; K=3; L=1; M=3.349417; N=10.048251; O=0.0602955800561437; P=0.0602955800561437; Q=815.8084429829045; R=7342.27598684614; S=0.04787034266582881; T=0.0478703426658288},
    @{A="FAPs"; B="Lama2"; C="Itga7"; D="FAPs"; E=3; F=1; G=243.5672963333334; H=730.701889; I=0.7939278902575405; J=0.7939278902575405; K=3; L=1; M=1.282876; N=3.848628; O=0.023094094452887; P=0.02309409445288699; Q=312.4666388509214; R=2812.199749658292; S=0.01833504568638894; T=0.01833504568638894},
    @{A="FAPs"; B="Lama2"; C="Itga7"; D="sCs"; E=3; F=1; G=243.5672963333334; H=730.701889; I=0.7939278902575405; J=0.7939278902575405; K=3; L=1; M=50.917666; N=152.752998; O=0.9166103254909692; P=0.9166103254909692; Q=12401.87824322369; R=111616.9041890132; S=0.7277225019053227; T=0.7277225019053227},
    @{A="sCs"; B="Lama2"; C="Itga7"; D="ECs"; E=3; F=1; G=60.76799533333334; H=182.303986; I=0.1980783424394845; J=0.1980783424394845; K=3; L=1; M=3.349417; N=10.048251; O=0.0602955800561437; P=0.0602955800561437; Q=203.5373566253874; R=1831.836209628486; S=0.01194324855394818; T=0.01194324855394818},
    @{A="sCs"; B="Lama2"; C="Itga7"; D="FAPs"; E=3; F=1; G=60.76799533333334; H=182.303986; I=0.1980783424394845; J=0.1980783424394845; K=3; L=1; M=1.282876; N=3.848628; O=0.023094094452887; P=0.02309409445288699; Q=77.95780278124535; R=701.6202250312081; S=0.00457443994936875; T=0.004574439949368749},
    @{A="sCs"; B="Lama2"; C="Itga7"; D="sCs"; E=3; F=1; G=60.76799533333334; H=182.303986; I=0.1980783424394845; J=0.1980783424394845; K=3; L=1; M=50.917666; N=152.752998; O=0.9166103254909692; P=0.9166103254909692; Q=3094.164489872225; R=27847.48040885003; S=0.1815606539361675; T=0.1815606539361675}
)

$cols = @("A","B","C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T")

$r = 2
foreach ($row in $newRows) {
    foreach ($c in $cols) {
        $ws.Range("$c$r").Value = $row[$c]
    }
    $r++
}
